# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Wed Dec 13 10:52:57 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.153.88'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '2.172.32'
$ws.Range('E3').Value = '  -2.36%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'251.25"
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('E6').Value = '  -3.07%  '
$ws.Range('D7').Value = "'65.91"
$ws.Range('E7').Value = '  -7.67%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = "'0.573"
$ws.Range('E9').Value = '  -4.83%  '
$ws.Range('D10').Value = "'58.96"
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('D11').Value = "'36.24"
$ws.Range('E11').Value = '  -11.60%  '
$ws.Range('D12').Value = "'0.0931"
$ws.Range('E12').Value = '  -3.81%  '
$ws.Range('D13').Value = "'0.103"
$ws.Range('E13').Value = '  -2.55%  '
$ws.Range('D14').Value = "'6.76"
$ws.Range('E14').Value = '  -6.37%  '
$ws.Range('D15').Value = '2.500.62'
$ws.Range('E15').Value = '  -2.17%  '
$ws.Range('D16').Value = "'14.22"
$ws.Range('E16').Value = '  -4.59%  '
$ws.Range('D17').Value = "'0.839"
$ws.Range('E17').Value = '  -2.76%  '
$ws.Range('D18').Value = '2.166.15'
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').Value = '41.036.62'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').Value = "'71.45"
$ws.Range('E21').Value = '  -1.99%  '
$ws.Range('D22').Value = "'6.03"
$ws.Range('E22').Value = '  -3.21%  '
$ws.Range('D23').Value = "'229.54"
$ws.Range('E23').Value = '  -2.34%  '
$ws.Range('D24').Value = "'2.01"
$ws.Range('E24').Value = '  -4.11%  '
$ws.Range('D25').Value = "'3.81"
$ws.Range('E25').Value = '  -5.73%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  +5.50%  '
$ws.Range('D28').Value = "'2.41"
$ws.Range('E28').Value = '  -4.84%  '
$ws.Range('D29').Value = "'2.12"
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').Value = "'168.74"
$ws.Range('E30').Value = '  -1.47%  '
$ws.Range('D31').Value = "'20.12"
$ws.Range('E31').Value = '  -3.21%  '
$ws.Range('E32').Value = '  -3.29%  '
$ws.Range('D33').Value = "'5.67"
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D34').Value = "'0.0745"
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('E35').Value = '  -3.12%  '
$ws.Range('D36').Value = "'4.49"
$ws.Range('E36').Value = '  -4.98%  '
$ws.Range('D37').Value = "'3.91"
$ws.Range('E37').Value = '  -2.42%  '
$ws.Range('D38').Value = "'24.37"
$ws.Range('E38').Value = '  -5.22%  '
$ws.Range('D39').Value = "'0.0303"
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = "'5.51"
$ws.Range('E40').Value = '  +13.10%  '
$ws.Range('E41').Value = '  -3.94%  '
$ws.Range('E42').Value = '  -8.31%  '
$ws.Range('D43').Value = "'60.79"
$ws.Range('E43').Value = '  -8.56%  '
$ws.Range('D44').Value = "'11.32"
$ws.Range('E44').Value = '  -7.30%  '
$ws.Range('D45').Value = "'8.44"
$ws.Range('E45').Value = '  -3.68%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = "'0.0987"
$ws.Range('E47').Value = '  -3.40%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = "'0.188"
$ws.Range('E48').Value = '  -7.76%  '
$ws.Range('D49').Value = "'1.13"
$ws.Range('E49').Value = '  -3.28%  '
$ws.Range('D50').Value = "'4.21"
$ws.Range('E50').Value = '  -9.55%  '
$ws.Range('E51').Value = '  -4.30%  '
